$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.868.98"
$ws.Range("E2").Value = "  +3.63%  "
$ws.Range("D3").Value = "2.648.05"
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'113.84"
$ws.Range("E5").Value = "  +7.12%  "
$ws.Range("D6").Value = "'326.74"
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.553"
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("D10").Value = "'41.18"
$ws.Range("E10").Value = "  +5.97%  "
$ws.Range("D11").Value = "'20.14"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "'7.35"
$ws.Range("E14").Value = "  +4.10%  "
$ws.Range("D15").Value = "3.063.09"
$ws.Range("E15").Value = "  +5.99%  "
$ws.Range("D16").Value = "2.645.01"
$ws.Range("E16").Value = "  +5.78%  "
$ws.Range("D17").Value = "'0.871"
$ws.Range("E17").Value = "  +4.52%  "
$ws.Range("D18").Value = "49.787.96"
$ws.Range("E18").Value = "  +3.81%  "
$ws.Range("D19").Value = "'13.16"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").Value = "0.0₃0955"
$ws.Range("E22").Value = "  +2.26%  "
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("D24").Value = "'277.07"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("D25").Value = "'2.58"
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("D26").Value = "'26.79"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "'10.04"
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("D30").Value = "'36.04"
$ws.Range("E30").Value = "  +3.59%  "
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "'50.37"
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("D33").Value = "'19.59"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("D35").Value = "'0.0806"
$ws.Range("E35").Value = "  +3.74%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  +6.74%  "
$ws.Range("D38").Value = "'4.75"
$ws.Range("E38").Value = "  +3.57%  "
$ws.Range("E39").Value = "  +6.53%  "
$ws.Range("D40").Value = "'125.43"
$ws.Range("E40").Value = "  +3.75%  "
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("D43").Value = "'22.20"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("E45").Value = "  +4.69%  "
$ws.Range("D46").Value = "2.073.20"
$ws.Range("E46").Value = "  +3.38%  "
$ws.Range("E47").Value = "  +13.28%  "
$ws.Range("D48").Value = "'1.98"
$ws.Range("E48").Value = "  +5.35%  "
$ws.Range("D49").Value = "'9.12"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("E50").Value = "  +4.40%  "
$ws.Range("D51").Value = "'81.51"
$ws.Range("E51").Value = "  +3.71%  "
